# Auto-generated edit script: updates crypto price (D) and volume change (E)
# columns for the cryptos.xlsx worksheet, matching the upstream data refresh
# commit "Updated cryptos list ... with GitHub Actions".
#
# Price cells (column D) are text-formatted in the source data (e.g. thousands
# separated by dots like "58.810.81", or plain decimals like "7.55"). Because
# many of the new price strings are valid numeric literals, we temporarily
# force the cell's number format to Text ("@") before assigning the value so
# Excel does not silently convert the string into a floating point number.
# The number format/style is reset back to "Normal" immediately afterward so
# no stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.810.81"
$ws.Range("E2").Value = "  -3.80%  "
$ws.Range("D3").Value = "3.214.85"
$ws.Range("E3").Value = "  -4.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.51%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.216.24"
$ws.Range("E8").Value = "  -4.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.458"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.64%  "
$ws.Range("D13").Value = "3.771.24"
$ws.Range("E13").Value = "  -4.37%  "
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.07%  "
$ws.Range("D16").Value = "3.217.59"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000157"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.38%  "
$ws.Range("D18").Value = "58.892.99"
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "361.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.519"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.45%  "
$ws.Range("D26").Value = "3.355.87"
$ws.Range("E26").Value = "  -4.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.171"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").Value = "0.0₃0958"
$ws.Range("E28").Value = "  -11.90%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.90%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -7.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0704"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.43%  "
$ws.Range("D42").Value = "3.247.13"
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.715"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.90%  "
$ws.Range("E47").Value = "  -6.89%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "2.298.21"
$ws.Range("E49").Value = "  -7.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.08%  "
